$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------
# Sheet1: BasicParsing - insert 5 front-matter rows at the top
# -------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("BasicParsing")
$ws1.Rows("1:5").Insert()

# Row 1: "---" front-matter delimiter (quote-prefixed text, left aligned,
# text number format - matches column A's existing text style)
$ws1.Range("A1").NumberFormat = "@"
$ws1.Range("A1").Value = "'---"

# Row 2: name / string / test basic spreadsheet parsing
$ws1.Range("A2:C4").HorizontalAlignment = -4131
$ws1.Range("A2").Value = "name"
$ws1.Range("B2").Value = "string"
$ws1.Range("C2").Value = "test basic spreadsheet parsing"

# Row 3: someNumber / number / 100
$ws1.Range("A3").Value = "someNumber"
$ws1.Range("B3").Value = "number"
$ws1.Range("C3").Value = 100

# Row 4: someBool / boolean / TRUE
$ws1.Range("A4").Value = "someBool"
$ws1.Range("B4").Value = "boolean"
$ws1.Range("C4").Value = $true

# Row 5: "---" front-matter delimiter
$ws1.Range("A5").NumberFormat = "@"
$ws1.Range("A5").Value = "'---"

# Data validation for the new numeric/bool column also gets DataTypeList
$ws1.Range("B2:B4").Validation.Add(3, 1, 1, "DataTypeList")

# -------------------------------------------------------------------
# Sheet2: FormulaAndRefParsing - just move the selection
# -------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("FormulaAndRefParsing")
$ws2.Range("E21").Select()

# -------------------------------------------------------------------
# Sheet3: ErrorCasesParsing - insert 3 front-matter rows at the top
# -------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("ErrorCasesParsing")
$ws3.Rows("1:3").Insert()

# Row 1: "---" front-matter delimiter (quote-prefixed, no extra numfmt
# since column A has no special default style on this sheet)
$ws3.Range("A1").Value = "'---"

# Row 2: created / date / 9/7/2024 (serial 45542)
$ws3.Range("A2").Value = "created"
$ws3.Range("B2").HorizontalAlignment = -4131
$ws3.Range("B2").Value = "date"
$ws3.Range("C2").NumberFormat = "mm-dd-yy"
$ws3.Range("C2").Value = 45542

# Row 3: "---" front-matter delimiter
$ws3.Range("A3").Value = "'---"

# Data validation for the new date cell also gets DataTypeList
$ws3.Range("B2").Validation.Add(3, 1, 1, "DataTypeList")

# -------------------------------------------------------------------
# Sheet4: .NestedDataParsing - just move the selection
# -------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(".NestedDataParsing")
$ws4.Range("H31").Select()

# -------------------------------------------------------------------
# Final view state: BasicParsing zoomed + selected at B3, then make
# ErrorCasesParsing the active/selected tab (matches activeTab=2)
# -------------------------------------------------------------------
$ws1.Range("B3").Select()
$excel.ActiveWindow.Zoom = 125

$ws3.Range("C3").Select()

Write-Output "done"
